$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Status text: a handback just completed, so every "Ready for handoff" cell
# (Overview + per-locale Status columns) now reads "Handed back: in sync
# with en-US".
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn: report the newly generated handback, by filling in the
# "Latest Target File" (F) / "Latest Handback File" (G) columns and
# refreshing the "Latest Handback DateTime" (H) column.
# ---------------------------------------------------------------------------
$aMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/aac54c978345d408a569a175111b0107a89f241c/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6280c72d048025a3d01ab963f3427d2559043cdf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)

$wsZhCn.Range("H2").Value = "2016-03-20 10:25:23"
$wsZhCn.Range("H3").Value = "2016-03-20 10:25:23"

# ---------------------------------------------------------------------------
# de-de: same handback report, plus its own (later) handback timestamp.
# ---------------------------------------------------------------------------
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b40576251e2db7e9cd34d58ece314fcc16de61c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)

$wsDeDe.Range("H2").Value = "2016-03-20 10:25:29"
$wsDeDe.Range("H3").Value = "2016-03-20 10:25:29"
